$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5: new customer record (customer4 / wheels / ajay / 72993 / "" / for fluid machinsm)
$ws.Range("A5").Value = "customer4"
$ws.Range("B5").Value = "wheels"
$ws.Range("C5").Value = "ajay"

# D5 must stay text ("72993"), not be auto-coerced into a number like Excel
# normally would do for a plain Value assignment. Build it as a text formula
# in a scratch cell, copy it, and paste-special as values onto D5 so the
# stored cell type ends up as a shared string instead of a numeric literal.
$ws.Range("H1").Formula = "=""72993"""
$ws.Range("H1").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("H1").ClearContents()

# E5 is blank in the source data.
$ws.Range("E5").Value = ""

$ws.Range("F5").Value = "for fluid machinsm"
